$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report date range) ---
$ws.Range("A8").Value = "Volume 31   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/4/2024  Through  11/10/2024"

# --- Weekly crime statistics table (rows 14-30) ---
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Range("C14").Value = 2
$ws.Range("C14").NumberFormat = "#,##0"
$ws.Range("F14").Value = 4
$ws.Range("I14").Value = 13
$ws.Range("K14").Value = 62.5
$ws.Range("L14").Value = -7.142857142857
$ws.Range("M14").Value = -7.142857142857
$ws.Range("N14").Value = -79.365079365079

# Row 15
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 38
$ws.Range("J15").Value = 38
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 35.714285714285
$ws.Range("M15").Value = 111.111111111111
$ws.Range("N15").Value = -47.945205479452

# Row 16
$ws.Range("C16").Value = 12
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 100
$ws.Range("G16").Value = 42
$ws.Range("H16").Value = 9.523809523809
$ws.Range("I16").Value = 552
$ws.Range("J16").Value = 580
$ws.Range("K16").Value = -4.827586206896
$ws.Range("L16").Value = 14.049586776859
$ws.Range("M16").Value = 40.458015267175
$ws.Range("N16").Value = -66.525166767738

# Row 17
$ws.Range("C17").Value = 14
$ws.Range("D17").Value = 18
$ws.Range("E17").Value = -22.222222222222
$ws.Range("F17").Value = 78
$ws.Range("G17").Value = 78
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 890
$ws.Range("J17").Value = 879
$ws.Range("K17").Value = 1.251422070534
$ws.Range("L17").Value = 20.759837177747
$ws.Range("M17").Value = 135.449735449735
$ws.Range("N17").Value = -8.530318602261

# Row 18
$ws.Range("C18").Value = 7
$ws.Range("E18").Value = 16.666666666666
$ws.Range("F18").Value = 18
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = -5.263157894736
$ws.Range("I18").Value = 299
$ws.Range("J18").Value = 274
$ws.Range("K18").Value = 9.12408759124
$ws.Range("L18").Value = -6.269592476489
$ws.Range("M18").Value = 70.857142857142
$ws.Range("N18").Value = -74.422583404619

# Row 19
$ws.Range("C19").Value = 19
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 90
$ws.Range("F19").Value = 83
$ws.Range("G19").Value = 70
$ws.Range("H19").Value = 18.571428571428
$ws.Range("I19").Value = 851
$ws.Range("J19").Value = 639
$ws.Range("K19").Value = 33.176838810641
$ws.Range("L19").Value = 31.733746130031
$ws.Range("M19").Value = 126.933333333333
$ws.Range("N19").Value = 27.014925373134

# Row 20
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 20
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 206
$ws.Range("J20").Value = 301
$ws.Range("K20").Value = -31.561461794019
$ws.Range("L20").Value = -23.703703703703
$ws.Range("M20").Value = 92.523364485981
$ws.Range("N20").Value = -64.543889845094

# Row 21
$ws.Range("D21").Value = 44
$ws.Range("E21").Value = 36.363636363636
$ws.Range("F21").Value = 250
$ws.Range("G21").Value = 231
$ws.Range("H21").Value = 8.225108225108
$ws.Range("I21").Value = 2849
$ws.Range("J21").Value = 2719
$ws.Range("K21").Value = 4.781169547627
$ws.Range("L21").Value = 14.051240992794
$ws.Range("M21").Value = 95.136986301369
$ws.Range("N21").Value = -44.978756276554

# Row 22
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 9
$ws.Range("H22").Value = 12.5
$ws.Range("I22").Value = 75
$ws.Range("J22").Value = 59
$ws.Range("K22").Value = 27.118644067796
$ws.Range("L22").Value = -3.846153846153
$ws.Range("M22").Value = 44.230769230769

# Row 23
$ws.Range("C23").Value = 9
$ws.Range("E23").Value = 12.5
$ws.Range("F23").Value = 31
$ws.Range("G23").Value = 37
$ws.Range("H23").Value = -16.216216216216
$ws.Range("I23").Value = 424
$ws.Range("J23").Value = 421
$ws.Range("K23").Value = 0.712589073634
$ws.Range("L23").Value = 24.340175953079
$ws.Range("M23").Value = 62.452107279693

# Row 24
$ws.Range("C24").Value = 52
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = 73.333333333333
$ws.Range("F24").Value = 202
$ws.Range("G24").Value = 117
$ws.Range("H24").Value = 72.649572649572
$ws.Range("I24").Value = 1503
$ws.Range("J24").Value = 1410
$ws.Range("K24").Value = 6.595744680851
$ws.Range("L24").Value = -10.429082240762
$ws.Range("M24").Value = 23.907666941467

# Row 25
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 17
$ws.Range("E25").Value = 29.411764705882
$ws.Range("F25").Value = 94
$ws.Range("G25").Value = 57
$ws.Range("H25").Value = 64.912280701754
$ws.Range("I25").Value = 654
$ws.Range("J25").Value = 540
$ws.Range("K25").Value = 21.111111111111
$ws.Range("L25").Value = -21.676646706586

# Row 26
$ws.Range("C26").Value = 27
$ws.Range("D26").Value = 26
$ws.Range("E26").Value = 3.846153846153
$ws.Range("F26").Value = 102
$ws.Range("G26").Value = 100
$ws.Range("H26").Value = 2
$ws.Range("I26").Value = 1097
$ws.Range("J26").Value = 981
$ws.Range("K26").Value = 11.824668705402
$ws.Range("L26").Value = 20.948180815876
$ws.Range("M26").Value = 9.371884346959

# Row 27
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -80
$ws.Range("I27").Value = 50
$ws.Range("J27").Value = 55
$ws.Range("K27").Value = -9.090909090909
$ws.Range("L27").Value = 8.695652173913

# Row 28
$ws.Range("C28").Value = 4
$ws.Range("F28").Value = 18
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 350
$ws.Range("I28").Value = 147
$ws.Range("K28").Value = 48.484848484848
$ws.Range("L28").Value = 93.421052631578

# Row 29
$ws.Range("C29").Value = 4
$ws.Range("D29").Value = 2
$ws.Range("D29").NumberFormat = "#,##0"
$ws.Range("E29").Value = 100
$ws.Range("E29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F29").Value = 8
$ws.Range("G29").Value = 6
$ws.Range("H29").Value = 33.333333333333
$ws.Range("I29").Value = 48
$ws.Range("J29").Value = 32
$ws.Range("K29").Value = 50
$ws.Range("L29").Value = -11.111111111111
$ws.Range("M29").Value = -11.111111111111
$ws.Range("N29").Value = -75.510204081632

# Row 30
$ws.Range("C30").Value = 4
$ws.Range("D30").Value = 2
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("E30").Value = 100
$ws.Range("E30").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F30").Value = 8
$ws.Range("G30").Value = 6
$ws.Range("H30").Value = 33.333333333333
$ws.Range("I30").Value = 43
$ws.Range("J30").Value = 30
$ws.Range("K30").Value = 43.333333333333
$ws.Range("L30").Value = -4.444444444444
$ws.Range("M30").Value = -6.521739130434
$ws.Range("N30").Value = -75.568181818181
